# Adds the 2024/11/21 column (BV) to the daily data table.
#
# Source data: a new date header in row 1 plus 52 numeric observations in
# rows 2-53, each carrying one of the sheet's three existing cell styles
# (s1 = no fill, s2 = yellow fill, s3 = light-blue fill) that are already
# used throughout the table to flag low-value readings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newCol = 74   # column BV -- the first column after the existing data (A:BU = 1:73)

# Exemplar cells (existing column BU = 73) to copy each named style from.
$styleSources = @{
    "s1" = $ws.Cells.Item(2, 73)   # BU2  -> s="1" (no fill)
    "s2" = $ws.Cells.Item(8, 73)   # BU8  -> s="2" (yellow fill)
    "s3" = $ws.Cells.Item(7, 73)   # BU7  -> s="3" (light-blue fill)
}

$rows = @(
    @{ Row=1; Style="s1"; Value="2024/11/21" }
    @{ Row=2; Style="s1"; Value="203.5" }
    @{ Row=3; Style="s3"; Value="128.1" }
    @{ Row=4; Style="s1"; Value="199.5" }
    @{ Row=5; Style="s3"; Value="125.9" }
    @{ Row=6; Style="s3"; Value="138" }
    @{ Row=7; Style="s1"; Value="345.2" }
    @{ Row=8; Style="s1"; Value="303.6" }
    @{ Row=9; Style="s1"; Value="183.6" }
    @{ Row=10; Style="s1"; Value="194.4" }
    @{ Row=11; Style="s2"; Value="124.3" }
    @{ Row=12; Style="s1"; Value="141.5" }
    @{ Row=13; Style="s1"; Value="163.8" }
    @{ Row=14; Style="s1"; Value="277.5" }
    @{ Row=15; Style="s3"; Value="136.6" }
    @{ Row=16; Style="s1"; Value="141.5" }
    @{ Row=17; Style="s1"; Value="154.7" }
    @{ Row=18; Style="s1"; Value="167.1" }
    @{ Row=19; Style="s1"; Value="211.7" }
    @{ Row=20; Style="s3"; Value="132.9" }
    @{ Row=21; Style="s1"; Value="223.3" }
    @{ Row=22; Style="s2"; Value="121.8" }
    @{ Row=23; Style="s3"; Value="137" }
    @{ Row=24; Style="s1"; Value="170.3" }
    @{ Row=25; Style="s1"; Value="151.2" }
    @{ Row=26; Style="s1"; Value="171.4" }
    @{ Row=27; Style="s1"; Value="185.4" }
    @{ Row=28; Style="s1"; Value="151.4" }
    @{ Row=29; Style="s1"; Value="183.8" }
    @{ Row=30; Style="s1"; Value="213.2" }
    @{ Row=31; Style="s3"; Value="128.1" }
    @{ Row=32; Style="s3"; Value="126.2" }
    @{ Row=33; Style="s1"; Value="204.4" }
    @{ Row=34; Style="s1"; Value="149.5" }
    @{ Row=35; Style="s3"; Value="138.9" }
    @{ Row=36; Style="s1"; Value="179.1" }
    @{ Row=37; Style="s1"; Value="163" }
    @{ Row=38; Style="s1"; Value="223.9" }
    @{ Row=39; Style="s1"; Value="197.2" }
    @{ Row=40; Style="s1"; Value="184.2" }
    @{ Row=41; Style="s3"; Value="136.8" }
    @{ Row=42; Style="s1"; Value="179.1" }
    @{ Row=43; Style="s1"; Value="148.1" }
    @{ Row=44; Style="s2"; Value="107.7" }
    @{ Row=45; Style="s1"; Value="156.7" }
    @{ Row=46; Style="s1"; Value="170.3" }
    @{ Row=47; Style="s1"; Value="184.7" }
    @{ Row=48; Style="s1"; Value="189.1" }
    @{ Row=49; Style="s3"; Value="135.2" }
    @{ Row=50; Style="s1"; Value="198.9" }
    @{ Row=51; Style="s1"; Value="157.1" }
    @{ Row=52; Style="s1"; Value="199.7" }
    @{ Row=53; Style="s3"; Value="136.7" }
)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r.Row, $newCol)

    if ($r.Row -eq 1) {
        # Row 1 holds the date label as text ("2024/11/21"). Assigning a
        # date-shaped string straight to .Value2 gets auto-converted into a
        # date serial number, so instead we go through a TEXT() formula and
        # then freeze it to a static value -- that keeps it a plain string.
        $cell.Formula = '=TEXT("' + $r.Value + '","yyyy/mm/dd")'
        $cell.Copy()
        $cell.PasteSpecial(-4163) | Out-Null   # xlPasteValues
    } else {
        $cell.Value2 = [double]$r.Value
    }

    $styleSources[$r.Style].Copy()
    $cell.PasteSpecial(-4122) | Out-Null       # xlPasteFormats
}

$excel.CutCopyMode = 0

# Register column BV's width (stored width 12, same as the other data
# columns) so a <col min="74" max="74" .../> entry is written out.
$ws.Columns.Item($newCol).ColumnWidth = 11.17

Write-Output "Done adding column BV (2024/11/21 data)"
